$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 3.075165666666667
$ws.Range("H2").Value2 = 9.225497
$ws.Range("I2").Value2 = 0.02641273658732285
$ws.Range("J2").Value2 = 0.02641273658732285
$ws.Range("M2").Value2 = 77.08952333333333
$ws.Range("N2").Value2 = 231.26857
$ws.Range("O2").Value2 = 0.2403816673726824
$ws.Range("P2").Value2 = 0.2403816673726824
$ws.Range("Q2").Value2 = 237.0630554143656
$ws.Range("R2").Value2 = 2133.56749872929
$ws.Range("S2").Value2 = 0.006349137660736119
$ws.Range("T2").Value2 = 0.006349137660736119

$ws.Range("G3").Value2 = 3.075165666666667
$ws.Range("H3").Value2 = 9.225497
$ws.Range("I3").Value2 = 0.02641273658732285
$ws.Range("J3").Value2 = 0.02641273658732285
$ws.Range("O3").Value2 = 0.3167483425780597
$ws.Range("P3").Value2 = 0.3167483425780597
$ws.Range("Q3").Value2 = 312.3754432261849
$ws.Range("R3").Value2 = 2811.378989035664
$ws.Range("S3").Value2 = 0.00836619053698539
$ws.Range("T3").Value2 = 0.008366190536985387

$ws.Range("G4").Value2 = 3.075165666666667
$ws.Range("H4").Value2 = 9.225497
$ws.Range("I4").Value2 = 0.02641273658732285
$ws.Range("J4").Value2 = 0.02641273658732285
$ws.Range("O4").Value2 = 0.4428699900492579
$ws.Range("P4").Value2 = 0.4428699900492579
$ws.Range("Q4").Value2 = 436.7559063047663
$ws.Range("R4").Value2 = 3930.803156742897
$ws.Range("S4").Value2 = 0.01169740838960134
$ws.Range("T4").Value2 = 0.01169740838960134

$ws.Range("I5").Value2 = 0.549422396165273
$ws.Range("J5").Value2 = 0.5494223961652731
$ws.Range("M5").Value2 = 77.08952333333333
$ws.Range("N5").Value2 = 231.26857
$ws.Range("O5").Value2 = 0.2403816673726824
$ws.Range("P5").Value2 = 0.2403816673726824
$ws.Range("Q5").Value2 = 4931.24790448771
$ws.Range("R5").Value2 = 44381.23114038939
$ws.Range("S5").Value2 = 0.1320710716821028
$ws.Range("T5").Value2 = 0.1320710716821028

$ws.Range("I6").Value2 = 0.549422396165273
$ws.Range("J6").Value2 = 0.5494223961652731
$ws.Range("O6").Value2 = 0.3167483425780597
$ws.Range("P6").Value2 = 0.3167483425780597
$ws.Range("S6").Value2 = 0.1740286333606163
$ws.Range("T6").Value2 = 0.1740286333606164

$ws.Range("I7").Value2 = 0.549422396165273
$ws.Range("J7").Value2 = 0.5494223961652731
$ws.Range("O7").Value2 = 0.4428699900492579
$ws.Range("P7").Value2 = 0.4428699900492579
$ws.Range("S7").Value2 = 0.2433226911225539
$ws.Range("T7").Value2 = 0.2433226911225539

$ws.Range("I8").Value2 = 0.424164867247404
$ws.Range("J8").Value2 = 0.4241648672474041
$ws.Range("M8").Value2 = 77.08952333333333
$ws.Range("N8").Value2 = 231.26857
$ws.Range("O8").Value2 = 0.2403816673726824
$ws.Range("P8").Value2 = 0.2403816673726824
$ws.Range("Q8").Value2 = 3807.020113067744
$ws.Range("R8").Value2 = 34263.1810176097
$ws.Range("S8").Value2 = 0.1019614580298435
$ws.Range("T8").Value2 = 0.1019614580298435

$ws.Range("I9").Value2 = 0.424164867247404
$ws.Range("J9").Value2 = 0.4241648672474041
$ws.Range("O9").Value2 = 0.3167483425780597
$ws.Range("P9").Value2 = 0.3167483425780597
$ws.Range("S9").Value2 = 0.134353518680458
$ws.Range("T9").Value2 = 0.1343535186804579

$ws.Range("I10").Value2 = 0.424164867247404
$ws.Range("J10").Value2 = 0.4241648672474041
$ws.Range("O10").Value2 = 0.4428699900492579
$ws.Range("P10").Value2 = 0.4428699900492579
$ws.Range("S10").Value2 = 0.1878498905371026
$ws.Range("T10").Value2 = 0.1878498905371027
